# Update gh-pages to output generated at 456a3b4
# Applies updated "want to go" counts (column F) and ticket-status text
# (column G) across the four worksheets of the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: 展览 (Exhibition)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2329
$ws.Range("F3").Value = 490
$ws.Range("F4").Value = 199
$ws.Range("F5").Value = 332
$ws.Range("F6").Value = 332
$ws.Range("F7").Value = 541
$ws.Range("F9").Value = 759
$ws.Range("F11").Value = 770
$ws.Range("F13").Value = 87
$ws.Range("F14").Value = 388
$ws.Range("F15").Value = 16
$ws.Range("F17").Value = 19804
$ws.Range("G17").Value = "已售罄"
$ws.Range("F18").Value = 600
$ws.Range("F19").Value = 63
$ws.Range("F20").Value = 221
$ws.Range("F22").Value = 170
$ws.Range("F23").Value = 142
$ws.Range("F26").Value = 205
$ws.Range("F29").Value = 136

# ---------------------------------------------------------------------
# Sheet: 演出 (Performance)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 177
$ws.Range("F6").Value = 202
$ws.Range("F8").Value = 3388
$ws.Range("F14").Value = 121
$ws.Range("F16").Value = 3053

# ---------------------------------------------------------------------
# Sheet: 本地生活 (Local Life)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 92
$ws.Range("F4").Value = 564

# ---------------------------------------------------------------------
# Sheet: 全部类型 (All Types)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 92
$ws.Range("F5").Value = 2329
$ws.Range("F6").Value = 564
$ws.Range("F7").Value = 490
$ws.Range("F8").Value = 199
$ws.Range("F9").Value = 332
$ws.Range("F10").Value = 332
$ws.Range("F11").Value = 541
$ws.Range("F12").Value = 177
$ws.Range("F16").Value = 202
$ws.Range("F18").Value = 759
$ws.Range("F20").Value = 770
$ws.Range("F22").Value = 87
$ws.Range("F23").Value = 388
$ws.Range("F24").Value = 16
$ws.Range("F26").Value = 19806
$ws.Range("G26").Value = "已售罄"
$ws.Range("F28").Value = 3388
$ws.Range("F32").Value = 600
$ws.Range("F33").Value = 63
$ws.Range("F34").Value = 221
$ws.Range("F38").Value = 170
$ws.Range("F39").Value = 142
$ws.Range("F42").Value = 121
$ws.Range("F44").Value = 205
$ws.Range("F47").Value = 136
$ws.Range("F48").Value = 3054
